$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the header (row 2), shifting existing
# data rows down by one, to add the latest day's price data.
$ws.Rows.Item(2).Insert()

# Force column A to be treated as plain text so the date-like string
# "2026-02-22" is not auto-converted into a date serial number, then
# restore the cell formatting back to the sheet's default afterward.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-22"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
